$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$src = $ws.Cells.Item(17,2)  # B17, style 2 in original (uses AGM154 weapons: F/A-18,F16,F14 -> shared string 41)
$dst = $ws.Cells.Item(19,2)
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
Write-Host "done"
